$d = $word.ActiveDocument

function Replace-Range([int]$start, [int]$len, [string]$newText) {
    $rng = $d.Range($start, $start + $len)
    $rng.Text = $newText
}

function Find-Offset([string]$needle) {
    $full = $d.Content.Text
    return $full.IndexOf($needle)
}

# ---------------------------------------------------------------
# Block 1: Title paragraph - "{{TIPO_PROJETO}} Nº {{NUMERO_PROJETO}}, DE {{DATA_PROJETO}}"
# Old runs: "{{TIPO_PROJETO}} Nº {{" | "NUMERO" | "_PROJETO}}, DE {{DATA_PROJETO}}"
# New runs: "{{TIPO_PROJETO}} Nº " | "{{NUMERO_PROJETO}}," | " DE {{DATA_PROJETO}}"
# ---------------------------------------------------------------
$idx = Find-Offset "{{TIPO_PROJETO}} N"
$oldFull = "{{TIPO_PROJETO}} Nº {{NUMERO_PROJETO}}, DE {{DATA_PROJETO}}"
Replace-Range $idx $oldFull.Length "##BLOCK1##"
Write-Output "Block1 merged"

# Now split the single placeholder run into three runs matching the target text.
# Replace placeholder (10 chars) with first chunk -> run becomes p1 only.
$idx = Find-Offset "##BLOCK1##"
$p1 = "{{TIPO_PROJETO}} Nº "
$p2 = "{{NUMERO_PROJETO}},"
$p3 = " DE {{DATA_PROJETO}}"
Replace-Range $idx 10 $p1
# Insert p2 as a new run right after p1 (zero-length range at the boundary creates a new run).
$idx2 = $idx + $p1.Length
Replace-Range $idx2 0 $p2
# Insert p3 as a new run right after p2.
$idx3 = $idx2 + $p2.Length
Replace-Range $idx3 0 $p3
Write-Output "Block1 split"
